$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force Text format on the Price/Volume columns so Excel does not
# auto-convert numeric-looking strings (e.g. "205.65") into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.923.58'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '1.545.59'
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '205.65'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").Value = '0.483'
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").Value = '1.766.12'
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").Value = '1.544.88'
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = '26.911.63'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '61.57'
$ws.Range("D18").Value = '213.55'
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").Value = '7.19'
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("E22").Value = '  -2.93%  '
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("E24").Value = '  -3.43%  '
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("D26").Value = '6.63'
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("E32").Value = '  +2.16%  '
$ws.Range("D33").Value = '1.365.79'
$ws.Range("E33").Value = '  -2.23%  '
$ws.Range("E34").Value = '  +0.70%  '
$ws.Range("D35").Value = '1.53'
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").Value = '0.968'
$ws.Range("E36").Value = '  +4.93%  '
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("E39").Value = '  -2.64%  '
$ws.Range("D40").Value = '0.804'
$ws.Range("E40").Value = '  -2.68%  '
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("D42").Value = '0.987'
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").Value = '5.49'
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("E44").Value = '  +1.63%  '
$ws.Range("D45").Value = '63.37'
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("E46").Value = '  -2.72%  '
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("D48").Value = '1.679.86'
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("D49").Value = '86.26'
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("E50").Value = '  +0.80%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '0.0949'
$ws.Range("E51").Value = '  -0.57%  '

# Restore the default (unstyled) look for the edited cells, matching the
# original workbook where these cells carried no explicit style index.
$ws.Range("D2:E51").Style = "Normal"
